$paragraphs = @(
    ,@(@{ text = "References:"; link = $null })
    ,@(@{ text = "https://github.com/peterhchen/300_Python_DataScience/upload/main"; link = "https://github.com/peterhchen/300_Python_DataScience/upload/main" })
    ,@(@{ text = "Python Interpreter and Environment:"; link = $null })
    ,@(@{ text = "1. "; link = $null }, @{ text = "https://www.python.org/"; link = "https://www.python.org/" })
    ,@(@{ text = "2. "; link = $null }, @{ text = "https://www.anaconda.com/"; link = "https://www.anaconda.com/" })
    ,@(@{ text = "3. "; link = $null }, @{ text = "https://www.jetbrains.com/pycharm/"; link = "https://www.jetbrains.com/pycharm/" })
    ,@(@{ text = "Languages, Data Structures, Networking, Multithreading, GUI Programming:"; link = $null })
    ,@(@{ text = "3. "; link = $null }, @{ text = "https://www.tutorialspoint.com/python/index.htm"; link = "https://www.tutorialspoint.com/python/index.htm" })
    ,@(@{ text = "2. "; link = $null }, @{ text = "https://www.w3schools.com/python/"; link = "https://www.w3schools.com/python/" })
    ,@(@{ text = "Numpy and Pandas:"; link = $null })
    ,@(@{ text = "1. "; link = $null }, @{ text = "https://numpy.org/"; link = "https://numpy.org/" })
    ,@(@{ text = "2. "; link = $null }, @{ text = "https://pandas.pydata.org/"; link = "https://pandas.pydata.org/" })
    ,@(@{ text = "PyTest"; link = $null }, @{ text = " and "; link = $null }, @{ text = "PyLint"; link = $null }, @{ text = ":"; link = $null })
    ,@(@{ text = "1. "; link = $null }, @{ text = "https://docs.pytest.org/en/6.2.x/"; link = "https://docs.pytest.org/en/6.2.x/" })
    ,@(@{ text = "2. "; link = $null }, @{ text = "https://www.pylint.org/"; link = "https://www.pylint.org/" })
    ,@(@{ text = "Textbook: “Python 3 for Absolute Beginners”, Time Hall and J-P Stacey, "; link = $null }, @{ text = "Apress"; link = $null }, @{ text = " (2000)"; link = $null })
    ,@(@{ text = "http://index-of.es/Python/Python%203%20for%20Absolute%20Beginners.pdf"; link = $null })
)

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$sh = $s.Shapes.Item(2)

# Build the full text (paragraphs joined with carriage returns) and set it in one shot
# so PowerPoint creates one <a:p> per paragraph, inheriting the bullet/run formatting
# already present on the placeholder's first paragraph.
$fullText = ""
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    $para = $paragraphs[$i]
    $paraText = ""
    foreach ($run in $para) {
        $paraText += $run.text
    }
    if ($i -eq 0) {
        $fullText = $paraText
    } else {
        $fullText += [char]13 + $paraText
    }
}

$tr = $sh.TextFrame.TextRange
$tr.Text = $fullText

# Now walk the paragraphs again, applying hyperlinks to the runs that need them.
$offset = 1
for ($i = 0; $i -lt $paragraphs.Count; $i++) {
    $para = $paragraphs[$i]
    foreach ($run in $para) {
        $len = $run.text.Length
        if ($run.link -ne $null) {
            $rng = $tr.Characters($offset, $len)
            $rng.ActionSettings(1).Hyperlink.Address = $run.link
        }
        $offset += $len
    }
    # account for the paragraph-mark character between paragraphs
    $offset += 1
}

# Resize the placeholder to match the new content height
$sh.Height = 391.22456692913386
